$d = $word.ActiveDocument

# Locate the paragraph containing "Step 1 (Instalasi MongoDB" and the
# "_GoBack" bookmark that currently sits at its start.
$stepPara = $d.Paragraphs.Item(4)
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Split a new paragraph in after it.
$stepPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(5)

# Insert the new paragraph's text, with a temporary trailing marker
# character so the bookmark's insertion point is mid-paragraph (not a
# paragraph-boundary position, which Word would otherwise push into the
# following paragraph).
$newPara.Range.InsertBefore("Membuka situs utamaX")
$newPara = $d.Paragraphs.Item(5)

$bookmarkPos = $newPara.Range.End - 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary marker character, leaving the bookmark collapsed
# right after the paragraph's text.
$markerRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$markerRange.Delete()
